$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '69.341.73'
$ws.Cells.Item(2, 5).Value = '  -2.22%  '
$ws.Cells.Item(3, 4).Value = '3.655.52'
$ws.Cells.Item(3, 5).Value = '  -4.56%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '673.61'
$ws.Cells.Item(5, 5).Value = '  -4.43%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '160.48'
$ws.Cells.Item(6, 5).Value = '  -6.44%  '
$ws.Cells.Item(7, 4).Value = '3.656.61'
$ws.Cells.Item(7, 5).Value = '  -4.43%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '1.01'
$ws.Cells.Item(8, 5).Value = '  +0.60%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.486'
$ws.Cells.Item(9, 5).Value = '  -7.29%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.147'
$ws.Cells.Item(10, 5).Value = '  -9.11%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '7.06'
$ws.Cells.Item(11, 5).Value = '  -4.60%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.443'
$ws.Cells.Item(12, 5).Value = '  -3.19%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000226'
$ws.Cells.Item(13, 5).Value = '  -11.04%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '4.287.02'
$ws.Cells.Item(14, 5).Value = '  -4.23%  '
$ws.Cells.Item(15, 2).Value = 'Avalanche'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '32.65'
$ws.Cells.Item(15, 5).Value = '  -10.61%  '
$ws.Cells.Item(16, 4).Value = '3.695.51'
$ws.Cells.Item(16, 5).Value = '  -5.54%  '
$ws.Cells.Item(17, 4).Value = '68.770.85'
$ws.Cells.Item(17, 5).Value = '  -3.22%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.113'
$ws.Cells.Item(18, 5).Value = '  -1.68%  '
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '16.05'
$ws.Cells.Item(19, 5).Value = '  -7.34%  '
$ws.Cells.Item(20, 2).Value = 'Polkadot'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.43'
$ws.Cells.Item(20, 5).Value = '  -10.52%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '478.34'
$ws.Cells.Item(21, 5).Value = '  -3.39%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '9.62'
$ws.Cells.Item(22, 5).Value = '  -9.77%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.650'
$ws.Cells.Item(23, 5).Value = '  -11.26%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '80.61'
$ws.Cells.Item(24, 5).Value = '  -5.61%  '
$ws.Cells.Item(25, 2).Value = 'WrappedeETH'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(25, 4).Value = '3.843.28'
$ws.Cells.Item(25, 5).Value = '  -3.54%  '
$ws.Cells.Item(26, 2).Value = 'Dai'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.995'
$ws.Cells.Item(26, 5).Value = '  -0.49%  '
$ws.Cells.Item(27, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.40'
$ws.Cells.Item(27, 5).Value = '  -5.50%  '
$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.0000123'
$ws.Cells.Item(28, 5).Value = '  -15.15%  '
$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '9.13'
$ws.Cells.Item(29, 5).Value = '  -13.90%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '2.71'
$ws.Cells.Item(30, 5).Value = '  -12.44%  '
$ws.Cells.Item(31, 2).Value = 'Fetch.AI'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.76'
$ws.Cells.Item(31, 5).Value = '  -15.39%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.01'
$ws.Cells.Item(32, 5).Value = '  -9.61%  '
$ws.Cells.Item(33, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.00'
$ws.Cells.Item(33, 5).Value = '  +0.30%  '
$ws.Cells.Item(34, 2).Value = 'NEARProtocol'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '6.56'
$ws.Cells.Item(34, 5).Value = '  -11.55%  '
$ws.Cells.Item(35, 2).Value = 'EthereumClassic'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '26.34'
$ws.Cells.Item(35, 5).Value = '  -10.19%  '
$ws.Cells.Item(36, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(36, 4).Value = '3.627.97'
$ws.Cells.Item(36, 5).Value = '  -4.36%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.157'
$ws.Cells.Item(37, 5).Value = '  -10.71%  '
$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '8.38'
$ws.Cells.Item(38, 5).Value = '  -8.68%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '6.03'
$ws.Cells.Item(39, 5).Value = '  +0.97%  '
$ws.Cells.Item(40, 2).Value = 'USDe'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  -0.01%  '
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0896'
$ws.Cells.Item(41, 5).Value = '  -12.37%  '
$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.17'
$ws.Cells.Item(42, 5).Value = '  -6.75%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 5).Value = '  -0.19%  '
$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.928'
$ws.Cells.Item(44, 5).Value = '  -11.21%  '
$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '48.07'
$ws.Cells.Item(45, 5).Value = '  -1.58%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '155.89'
$ws.Cells.Item(46, 5).Value = '  -4.39%  '
$ws.Cells.Item(47, 2).Value = 'dogwifhat'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.77'
$ws.Cells.Item(47, 5).Value = '  -16.22%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.000271'
$ws.Cells.Item(48, 5).Value = '  -13.10%  '
$ws.Cells.Item(49, 2).Value = 'ONDO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.27'
$ws.Cells.Item(49, 5).Value = '  -7.31%  '
$ws.Cells.Item(50, 2).Value = 'Cosmos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '7.86'
$ws.Cells.Item(50, 5).Value = '  -10.06%  '
$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '27.31'
$ws.Cells.Item(51, 5).Value = '  -3.27%  '
